$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

# Title (appears twice: Heading1 and bold paragraph near the end)
Replace-Text "Play Balloonies for Free: Review and Game Overview | IGT" "Play Balloonies Free - Exciting Gameplay and Fantastic Bonuses"

# "What we like" bullet list
Replace-Text "Refreshing and unique gameplay." "Refreshing and unique gameplay dynamic"
Replace-Text "Excellent graphic quality with colorful and cartoonish design." "Colorful and cartoonish design with excellent graphics"
Replace-Text "Exciting array of bonus features." "Exciting array of bonus features"
Replace-Text "Decent chance of winning with multipliers and modifiers." "Decent chance of winning with good payouts"

# "What we don't like" bullet list
Replace-Text "Low to medium payouts." "Limited availability of the bonus features"
Replace-Text "Limited number of free spins." "RTP could be higher"

# Meta description (italic paragraph at the end)
Replace-Text "Learn about the exciting bonus features and unique gameplay of Balloonies by IGT, and play it for free before trying your luck with real money." "Play Balloonies for free and enjoy the unique gameplay. Experience exciting bonuses and good chances of winning."
